# Swap the two embedded themes in the deck:
#   ppt/theme/theme1.xml  "Integral"      -> "Office Theme" colours
#   ppt/theme/theme2.xml  "Office Theme"  -> "Integral" colours
#
# The font scheme / format scheme are already byte-identical between the
# two themes, so the only observable difference is the colour scheme
# (and its `name` attribute, which PowerPoint's object model does not
# expose for writing). We drive the swap through the slide master's
# Theme.ThemeColorScheme, which is what the OOXML <a:clrScheme> maps to.

function Convert-HexToPptRgb {
    param([string]$Hex)
    $r = [Convert]::ToInt32($Hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($Hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($Hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Target "Office Theme" colour scheme (was theme2.xml, now becomes theme1.xml)
$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$master = $p.SlideMaster
$theme = $master.Theme
$scheme = $theme.ThemeColorScheme

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $scheme.Item($i).RGB = Convert-HexToPptRgb $officeColors[$i - 1]
}

# Best-effort: reflect the new theme's display name too.
$theme.Name = "Office Theme"
